$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header labels
$ws.Range("B2").Value = "Locomotion"

# Add new sub-header in G3 (added before the other new strings, to match shared-string order)
$ws.Range("G3").Value = "Wheel Diameter (cm)"

# Add new "Component Sizes" header in G2
$ws.Range("G2").Value = "Component Sizes"

# Update Bounding Box units label last
$ws.Range("D2").Value = "Bounding Box (cm)"

# Add new data row 4 with values
$ws.Range("B4").Value = 4.5
$ws.Range("C4").Value = 0.2
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 9

# Set column G width as close as possible to the target stored width (19.5703125 characters).
# Note: this runtime's COM ColumnWidth setter snaps the persisted width to 1/6-character
# increments, so we pick the input that lands on the closest achievable value (19.5).
$ws.Columns.Item(7).ColumnWidth = 18.6667

# Update selection to C24 as in the target file
$ws.Range("C24").Select()
